$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Weight_before_[mg] / Weight_after_[mg]
# shift right to E and F respectively) and populate it with the bamboo species.
$ws.Columns("D").Insert()

$ws.Range("D1").Value = "Bamboo_sp"

$ws.Range("D2").Value = "Bambusa blumeana"
$ws.Range("D3").Value = "Bambusa blumeana"
$ws.Range("D4").Value = "Bambusa blumeana"
$ws.Range("D5").Value = "Bambusa blumeana"
$ws.Range("D6").Value = "Bambusa blumeana"
$ws.Range("D7").Value = "Bambusa blumeana"

$ws.Range("D8").Value = "Schizostachum lima"
$ws.Range("D9").Value = "Schizostachum lima"
$ws.Range("D10").Value = "Schizostachum lima"
$ws.Range("D11").Value = "Schizostachum lima"
$ws.Range("D12").Value = "Schizostachum lima"
$ws.Range("D13").Value = "Schizostachum lima"
